$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row appended after the 2025-10-26 entry, for the 2025-11-02 run.
# The date column stores the date as plain text (not an Excel date
# serial), matching the other recently-appended rows in this sheet.
$ws.Cells.Item(55, 1).NumberFormat = "@"
$ws.Cells.Item(55, 1).Value = "11/02/2025"
$ws.Cells.Item(55, 1).Style = "Normal"

$ws.Cells.Item(55, 2).Value = 0.0004510500000000014
$ws.Cells.Item(55, 3).Value = 110852.4553818863
$ws.Cells.Item(55, 4).Value = 50
